$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Carry the bordered, no-fill row formatting down onto the two new rows
# (same visual style already used for every data row on this sheet).
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E26").PasteSpecial(-4122)

# Row 25 - TestCase_A24
$ws.Range("A25").Value = "TestCase_A24"
$ws.Range("C25").Value = "Verify that TR account gets locked after 5 consecutive unsuccessful login attempts"
$ws.Range("B25").Value = "OPQA-525"
$ws.Range("D25").Value = "Y"
$ws.Range("E25").Value = "SKIP"

# Row 26 - TestCase_A25
$ws.Range("A26").Value = "TestCase_A25"
$ws.Range("B26").Value = "OPQA-529"
$ws.Range("C26").Value = "Verify that Help link is working properly"
$ws.Range("D26").Value = "Y"
$ws.Range("E26").Value = "PASS"

$ws.Range("D19").Select()
